# This script applies a weekly data rotation to the Fruta/Hortalizas
# 'Granada' worksheet: the market-report values (date, variety, quality,
# volume, prices, unit, origin, etc.) held in rows 2,3,4,5,6,8,9,10,13
# are rotated one step along the cycle 2 -> 5 -> 13 -> 6 -> 10 -> 4 -> 9 -> 3 -> 8 -> 2
# (row 7, 11 and 12 are left untouched), matching the published weekly update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44334
$ws.Range("M2").Value = 16
$ws.Range("N2").Value = 240000
$ws.Range("O2").Value = 250000
$ws.Range("P2").Value = 245000
$ws.Range("Q2").Value = "`$/bins (450 kilos)"
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("S2").Value = 544
$ws.Range("T2").Value = 450

# Row 3
$ws.Range("D3").Value = 44312
$ws.Range("K3").Value = "Wonderfull"
$ws.Range("M3").Value = 24
$ws.Range("N3").Value = 220000
$ws.Range("O3").Value = 240000
$ws.Range("P3").Value = 230000
$ws.Range("Q3").Value = "`$/bins (400 kilos)"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 575
$ws.Range("T3").Value = 400

# Row 4
$ws.Range("D4").Value = 44312
$ws.Range("K4").Value = "Wonderfull"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 34
$ws.Range("N4").Value = 240000
$ws.Range("O4").Value = 240000
$ws.Range("P4").Value = 240000
$ws.Range("Q4").Value = "`$/bins (450 kilos)"
$ws.Range("R4").Value = "Provincia del Elquí"
$ws.Range("S4").Value = 533
$ws.Range("T4").Value = 450

# Row 5
$ws.Range("D5").Value = 44285
$ws.Range("M5").Value = 8
$ws.Range("N5").Value = 280000
$ws.Range("O5").Value = 300000
$ws.Range("P5").Value = 290000
$ws.Range("R5").Value = "Provincia del Elquí"
$ws.Range("S5").Value = 725

# Row 6
$ws.Range("D6").Value = 44721
$ws.Range("K6").Value = "Wonderfull"
$ws.Range("M6").Value = 7
$ws.Range("N6").Value = 300000
$ws.Range("O6").Value = 300000
$ws.Range("P6").Value = 300000
$ws.Range("Q6").Value = "`$/bins (400 kilos)"
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 750
$ws.Range("T6").Value = 400

# Row 8
$ws.Range("D8").Value = 44662
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("M8").Value = 45
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = "`$/caja 18 kilos granel"
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44662
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 16000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 16000
$ws.Range("Q9").Value = "`$/caja 18 kilos granel"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 889
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44307
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("M10").Value = 150
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 17000
$ws.Range("Q10").Value = "`$/caja 15 kilos granel"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1133
$ws.Range("T10").Value = 15

# Row 13
$ws.Range("D13").Value = 44320
$ws.Range("M13").Value = 12
$ws.Range("N13").Value = 250000
$ws.Range("O13").Value = 260000
$ws.Range("P13").Value = 255000
$ws.Range("R13").Value = "Provincia de Limarí"
$ws.Range("S13").Value = 638
